$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (K-Means++)
$ws.Range("B2").Value = 328.6670608596
$ws.Range("C2").Value = 13.8896666801
$ws.Range("D2").Value = 0.2009297849
$ws.Range("E2").Value = 0.0124132872
$ws.Range("F2").Value = 1.0194684532
$ws.Range("G2").Value = 0.0297859789

# Row 3 (PSO)
$ws.Range("B3").Value = 391.6973940503
$ws.Range("C3").Value = 36.1625748871
$ws.Range("D3").Value = 0.1343657739
$ws.Range("E3").Value = 0.0516522078
$ws.Range("F3").Value = 1.0123951504
$ws.Range("G3").Value = 0.0508387423

# Row 4 (PSO Hybrid)
$ws.Range("B4").Value = 323.396516429
$ws.Range("C4").Value = 3.0138386068
$ws.Range("D4").Value = 0.20871241
$ws.Range("E4").Value = 0.0094420556
$ws.Range("F4").Value = 1.0020165163
$ws.Range("G4").Value = 0.004504653
